# 9.5.1.xlsx — add a new "2023" data column (L) to the existing table.
#
# The sheet holds a small table with year headers in row 4 (D4:K4, ending
# at 2022 in column K) and the corresponding percentage values in row 5
# (D5:K5), plus a thin bottom-border spacer row above it (row 3). We extend
# that table one column to the right (column L) for the new 2023 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column L: replicate the existing formatting from column K (border
#     spacer row, year-header row, value row) before writing the new data,
#     so the new cells pick up the same styles (borders/number formats)
#     as their column-K neighbours. ---

$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L4").Value = 2023

$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L5").Value = 0.11972285283622097

$excel.CutCopyMode = 0

# Row 5 (the value row) grows a bit taller to fit the now-wider table.
$ws.Rows.Item(5).RowHeight = 40.5
